$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Input data added for the five test-case rows. Columns:
#   B: input: Time btw births (s)
#   C: input: Time btw deaths (s)
#   D: input: Time btw net migrations (s)
#   E: Current Population
#   F: No of Years in future projections
# Output columns (G:I -> Population Change / Future Population / Increase
# Decrease) are intentionally left blank, per the commit message
# ("input added but no output").
$rows = @(
    @{ Row = 6;  B = 8;  C = 12; D = 126; E = 333100360; F = 5 },
    @{ Row = 7;  B = 5;  C = 12; D = 126; E = 333100360; F = 5 },
    @{ Row = 8;  B = 8;  C = 12; D = 126; E = 333100360; F = 5 },
    @{ Row = 9;  B = 8;  C = 8;  D = 126; E = 333100360; F = 5 },
    @{ Row = 10; B = 14; C = 12; D = 200; E = 333100360; F = 5 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# Move the active selection, matching where the author's cursor ended up.
$ws.Range("I13").Select()
